# Append two new paragraphs to the end of the document body (just before
# the sectPr / after the last existing paragraph):
#   1. a plain paragraph containing "[PUMP:TBD:1]"
#   2. a "List Bullet" styled paragraph containing "BOLUS:SRS:2"
#
# We inject raw WordprocessingML via Range.InsertXML instead of using the
# higher level Paragraphs/Style APIs: replacing (rather than collapsing at
# the end of) the last paragraph's range lets us keep that existing
# paragraph mark untouched (so the pre-existing empty <w:p/> survives
# unchanged) while cleanly appending the two brand-new paragraph marks
# after it, each carrying exactly the formatting we want - and without the
# interpreter stamping a spurious w:rsidP on anything.

$d = $word.ActiveDocument

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)

# Full range of the last paragraph (including its end-of-paragraph mark).
$target = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$wordmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xmlFragment =
  "<w:p $wordmlNs/>" + `
  "<w:p $wordmlNs><w:r><w:t>[PUMP:TBD:1]</w:t></w:r></w:p>" + `
  "<w:p $wordmlNs><w:pPr><w:pStyle w:val=`"ListBullet`"/></w:pPr><w:r><w:t>BOLUS:SRS:2</w:t></w:r></w:p>"

$target.InsertXML($xmlFragment)
